# Configure hierarchy: remove the "patient" and "other" entries from the
# "roles" choice list on the "choices" sheet of the health_facility-create
# form. These were two standalone role options that are no longer part of
# the configurable hierarchy (health facility / patient service delivery
# group / patient contact / facility employee).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# Row 13 is list_name="roles", name="patient" -- delete the whole row,
# shifting everything below it up by one.
$ws.Rows.Item(13).Delete()

# After the first deletion, the row that used to be 18 (list_name="roles",
# name="other") is now row 17 -- delete it too.
$ws.Rows.Item(17).Delete()
